# Add new KEYWORDS / CATEGORY rows to the "Feed" worksheet (rows 144-166).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data to append: Column C = keyword, Column D = category
$newRows = @(
    @{ C = "satellite operators";              D = "Satellite Manufacturers" },
    @{ C = "satellite manufacturing";           D = "Satellite Manufacturers" },
    @{ C = "space component manufacturing";     D = "Satellite Manufacturers" },
    @{ C = "constellations";                    D = "Satellite Manufacturers" },
    @{ C = "optical communication ";            D = "Satellite Manufacturers" },
    @{ C = "IOD";                                D = "Satellite Manufacturers" },
    @{ C = "IOV";                                D = "Satellite Manufacturers" },
    @{ C = "Hosted payload";                     D = "Satellite Services" },
    @{ C = "contracts";                          D = "Finance & Deals" },
    @{ C = "APAC";                               D = "Finance & Deals" },
    @{ C = "finance";                            D = "Finance & Deals" },
    @{ C = "deal";                               D = "Finance & Deals" },
    @{ C = "funding";                            D = "Finance & Deals" },
    @{ C = "NAM";                                D = "Finance & Deals" },
    @{ C = "broadband";                          D = "Satellite Manufacturers" },
    @{ C = "launch";                             D = "Launches" },
    @{ C = "deployment ";                        D = "Satellite Operations" },
    @{ C = "filings";                            D = "Satellite Operations" },
    @{ C = "Amazon Leo ";                        D = "Satellite Operations" },
    @{ C = "SAR";                                D = "Satellite Manufacturers" },
    @{ C = "space agency";                       D = "Institutional" },
    @{ C = "space agencies";                     D = "Institutional" },
    @{ C = "institutional";                      D = "Institutional" }
)

$startRow = 144
$endRow = $startRow + $newRows.Count - 1

# Column C (keywords) is filled in first, then column D (categories) is
# filled in afterwards -- this mirrors how the shared strings table ends up
# ordered (all keyword strings first, the handful of category strings last).
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 3).Value = $newRows[$i].C
}
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 4).Value = $newRows[$i].D
}

# Match the row height used by the rest of the keyword table.
$ws.Range("A$startRow`:A$endRow").RowHeight = 15.75

# Update the view state to reflect the new selection / scroll position.
$lastCell = $ws.Range("D" + ($endRow + 1))
$lastCell.Select()
$excel.ActiveWindow.ScrollRow = 146
